$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$origStyle = $ws.Range("D2").Style

$ws.Range("D2").Value = "68.172.38"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "3.884.48"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "481.83"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.86"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -3.64%  "

$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.738"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +1.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +7.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000356"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.80"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.50"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").Value = "4.508.69"
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("D15").Value = "3.912.22"
$ws.Range("E15").Value = "  -1.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.21"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -4.15%  "

$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.90"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").Value = "68.208.33"
$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "428.85"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -2.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.76"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.44"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +1.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.97"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +13.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.66"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.36"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.65"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -3.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "709.90"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.51"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +1.29%  "

$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("E33").Value = "  +2.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.09"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +9.43%  "

$ws.Range("D35").Value = "0.0₃0880"
$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.85"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -2.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.75"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +3.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0503"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +6.61%  "

$ws.Range("E39").Value = "  -4.47%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.394"
$ws.Range("D41").Style = $origStyle

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.94"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +2.35%  "

$ws.Range("E43").Value = "  +2.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.97"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("E46").Value = "  +3.96%  "

$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.37"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -2.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -3.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.64"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -2.00%  "
